# Planilha1: Apache POI demo workbook - adding sum/multiplication formulas
# so that after the sheet is generated, the user can change the inputs and
# the formulas keep the dependent cells (including the Planilha2 totals)
# updated without errors.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Planilha1")

# D4: lucro base value changes from 100 to 520
$ws1.Range("D4").Value = 520

# Row 6: new sum/multiplication example (quantidade/valor -> soma/multiplicacao)
$ws1.Range("D6").Value = 500.3
$ws1.Range("E6").Value = 2
$ws1.Range("F6").Formula = "=SUM(D6:E6)"
$ws1.Range("G6").Formula = "=D6*E6"

# Row heights: rows 6 and 7 shrink slightly, and a new (empty) row 8 appears
$ws1.Rows.Item(6).RowHeight = 13.8
$ws1.Rows.Item(7).RowHeight = 13.8
$ws1.Rows.Item(8).RowHeight = 13.8

# Touch row 8 (formatting-only) so the worksheet's used range / dimension
# grows to include it, matching the new B2:G8 extent, without putting any
# real value into the row.
$ws1.Range("B8").NumberFormat = "General"

# Move the active selection to D5
$ws1.Range("D5").Select()

# Sheet2!D2 depends on Planilha1!D4 - C2 and recalculates automatically.
